{"js": "// 1) \"Headings Detected: 8\" -> \"Headings Detected: 5\"\nconst headingCountResults = context.document.body.search(\"Headings Detected: 8\", { matchCase: true });\nheadingCountResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < headingCountResults.items.length; i++) {\n  headingCountResults.items[i].insertText(\"Headings Detected: 5\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Strip the yellow-highlight + blue \"[HEADING Level 1]\" annotation runs\n//    from the three author/affiliation paragraphs, leaving plain text only.\nconst targetTexts = [\n  \"Rohit Kumar, Ananya Sharma, Vikram Patel\",\n  \"Department of Computer Science\",\n  \"XYZ University\"\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const fullText = para.text.replace(/\\r$/, \"\");\n  const matchedTarget = targetTexts.find(\n    (t) => fullText === t || fullText === t + \" [HEADING Level 1]\"\n  );\n  if (matchedTarget) {\n    const cleanOoxml =\n      '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p><w:r><w:t>' + escapeXml(matchedTarget) + '</w:t></w:r></w:p></w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData></pkg:part></pkg:package>';\n    para.insertOoxml(cleanOoxml, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Headings Detected: 8\" -> \"Headings Detected: 5\"\n$find = $d.Content.Find\n$find.Execute(\"Headings Detected: 8\", $false, $false, $false, $false, $false, $true, 1, $false, \"Headings Detected: 5\", 2) | Out-Null\n\n# 2) Strip the yellow-highlight run + the blue \"[HEADING Level 1]\" annotation\n#    run from the three author/affiliation paragraphs, leaving plain text.\n$targets = @(\n    \"Rohit Kumar, Ananya Sharma, Vikram Patel\",\n    \"Department of Computer Science\",\n    \"XYZ University\"\n)\n\nforeach ($t in $targets) {\n    $annotated = \"$t [HEADING Level 1]\"\n\n    # Merge the two runs back into a single plain run by replacing the\n    # annotated text (base text + \" [HEADING Level 1]\") with just the base\n    # text. Find/Replace keeps the first run's formatting (yellow highlight)\n    # and drops the trailing blue-colored annotation run entirely.\n    $find2 = $d.Content.Find\n    $found = $find2.Execute($annotated, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)\n\n    if ($found) {\n        # Now strip the leftover yellow highlight from that run so only\n        # plain, unformatted text remains, matching the target.\n        $find3 = $d.Content.Find\n        $found3 = $find3.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n        if ($found3) {\n            $find3.Parent.HighlightColorIndex = 0\n        }\n    }\n}\n\nWrite-Output \"done\"\n"}
